$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new string values in the order NoraUi's bakery sample introduces
# them, so the shared-string table is built up in the same first-seen order
# as the target workbook: masked password, sgrillon, sgrillon2, profile,
# admin, referencer. "Result" already exists in the shared string table and
# is simply reused for the new D1 header.
$ws.Range("B2").Value = "℗:qmTAYKS9UG87rNuUQ0Ao6Q=="
$ws.Range("A2").Value = "sgrillon"
$ws.Range("A3").Value = "sgrillon2"
$ws.Range("B3").Value = "℗:qmTAYKS9UG87rNuUQ0Ao6Q=="
$ws.Range("C1").Value = "profile"
$ws.Range("C2").Value = "admin"
$ws.Range("C3").Value = "referencer"
$ws.Range("D1").Value = "Result"

# D1 takes over the bold/filled header style that used to belong to C1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# C1 becomes a regular header cell, matching the plain A1/B1 style
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Column widths (target widths of 11.140625 / 34.140625 / 12.140625 / 8.42578125
# characters as stored in the OOXML; the engine quantizes ColumnWidth to whole
# pixels, so feed it the values that round-trip to the closest achievable width)
$ws.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws.Columns.Item(2).ColumnWidth = 33.333333333333336
$ws.Columns.Item(3).ColumnWidth = 11.333333333333334
$ws.Columns.Item(4).ColumnWidth = 7.666666666666667

# Update selection to match the target workbook
$ws.Range("B3").Select() | Out-Null
